$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url4 = "https://www.genomeweb.com/cancer/entrogen-colorectal-cancer-ras-mutation-detection-test-nabs-cms-coverage"
$url5 = "https://www.360dx.com/cancer/entrogen-colorectal-cancer-ras-mutation-detection-test-nabs-cms-coverage"
$title = "EntroGen Colorectal Cancer RAS Mutation Detection Test Nabs CMS Coverage"
$keyword = "companion diagnostic"

# Row 4
$ws.Range("A4").Value = $url4
$ws.Hyperlinks.Add($ws.Range("A4"), $url4)
$ws.Range("A4").Style = $ws.Range("A2").Style
$ws.Range("B4").Value = $keyword
$ws.Range("C4").Value = $title

# Row 5
$ws.Range("A5").Value = $url5
$ws.Hyperlinks.Add($ws.Range("A5"), $url5)
$ws.Range("A5").Style = $ws.Range("A2").Style
$ws.Range("B5").Value = $keyword
$ws.Range("C5").Value = $title
